$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Info")

# Wrap the @base / @prefix URI values (column D) in angle brackets.
$ws.Cells.Item(1, 4).Value2 = "<" + $ws.Cells.Item(1, 4).Value2 + ">"
$ws.Cells.Item(2, 4).Value2 = "<" + $ws.Cells.Item(2, 4).Value2 + ">"
$ws.Cells.Item(3, 4).Value2 = "<" + $ws.Cells.Item(3, 4).Value2 + ">"
$ws.Cells.Item(4, 4).Value2 = "<" + $ws.Cells.Item(4, 4).Value2 + ">"
$ws.Cells.Item(5, 4).Value2 = "<" + $ws.Cells.Item(5, 4).Value2 + ">"

# Widen column D to fit the new, longer values.
$ws.Range("D4:D4").EntireColumn.ColumnWidth = 29.140625

# Move the selection to D1:D5 (mirrors the user reviewing the freshly
# updated column D after wrapping the URIs in angle brackets).
$ws.Range("D1:D5").Select()
